$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2244.027
$ws.Range("J17").Value = 2244.027
$ws.Range("L17").Value = 6732.081
$ws.Range("N17").Value = -7068.081
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 1500
$ws.Range("M51").Value = -1016
$ws.Range("H82").Value = 17600
$ws.Range("I82").Value = 8000
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 24000
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -23594
$ws.Range("N82").Value = -60812
$ws.Range("H85").Value = 17600
$ws.Range("I85").Value = 8000
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 24000
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -22596
$ws.Range("N85").Value = -62808
$ws.Range("H125").Value = 420
$ws.Range("J125").Value = 260
$ws.Range("L125").Value = 2340
$ws.Range("N125").Value = -7260
$ws.Range("H132").Value = 5549.5
$ws.Range("I132").Value = 5554.6665
$ws.Range("J132").Value = 5544.3335
$ws.Range("K132").Value = 16663.9995
$ws.Range("L132").Value = 16633.0005
$ws.Range("M132").Value = -14133.9995
$ws.Range("N132").Value = -21693.0005
$ws.Range("H137").Value = 2480.6365
$ws.Range("I137").Value = 2976.7144
$ws.Range("J137").Value = 1612.5
$ws.Range("K137").Value = 8930.143199999999
$ws.Range("L137").Value = 4837.5
$ws.Range("M137").Value = -6380.143199999999
$ws.Range("N137").Value = -9937.5
$ws.Range("H138").Value = 4596.4443
$ws.Range("J138").Value = 4995.4
$ws.Range("L138").Value = 14986.2
$ws.Range("N138").Value = -25266.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9899.25
$ws.Range("J2").Value = 9899.5
$ws.Range("L2").Value = 9899.5
$ws.Range("N2").Value = -10125.5
$ws.Range("H5").Value = 96.75
$ws.Range("I5").Value = 99
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 99
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 13
$ws.Range("N5").Value = -314
$ws.Range("H45").Value = 2996.25
$ws.Range("I45").Value = 2998.3333
$ws.Range("K45").Value = 2998.3333
$ws.Range("M45").Value = -2621.3333
$ws.Range("H116").Value = 9899.25
$ws.Range("J116").Value = 9899.5
$ws.Range("L116").Value = 9899.5
$ws.Range("N116").Value = -14487.5
$ws.Range("H122").Value = 6812
$ws.Range("I122").Value = 6812
$ws.Range("K122").Value = 20436
$ws.Range("M122").Value = -17986

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9899.25
$ws.Range("J3").Value = 9899.5
$ws.Range("L3").Value = 9899.5
$ws.Range("N3").Value = -10127.5
$ws.Range("H4").Value = 96.75
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 16
$ws.Range("N4").Value = -320

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 311.8
$ws.Range("I22").Value = 314.75
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 314.75
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 35.25
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 791.3570999999999
$ws.Range("I31").Value = 839.9167
$ws.Range("K31").Value = 839.9167
$ws.Range("M31").Value = -544.9167
$ws.Range("H34").Value = 791.3570999999999
$ws.Range("I34").Value = 839.9167
$ws.Range("K34").Value = 839.9167
$ws.Range("M34").Value = -637.9167
$ws.Range("H58").Value = 2559.4
$ws.Range("I58").Value = 2715.3215
$ws.Range("J58").Value = 1935.7142
$ws.Range("K58").Value = 2715.3215
$ws.Range("L58").Value = 1935.7142
$ws.Range("M58").Value = -2512.3215
$ws.Range("N58").Value = -2341.7142
$ws.Range("H62").Value = 5201.6665
$ws.Range("I62").Value = 4802.5
$ws.Range("K62").Value = 4802.5
$ws.Range("M62").Value = -4178.5
$ws.Range("H65").Value = 5201.6665
$ws.Range("I65").Value = 4802.5
$ws.Range("K65").Value = 24012.5
$ws.Range("M65").Value = -20892.5
$ws.Range("H86").Value = 7999
$ws.Range("I86").Value = 7999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7999
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -6876
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 7999
$ws.Range("I89").Value = 7999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 39995
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -34379
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 6582.6665
$ws.Range("I134").Value = 7099.4
$ws.Range("K134").Value = 21298.2
$ws.Range("M134").Value = -18763.2
$ws.Range("H136").Value = 2559.4
$ws.Range("I136").Value = 2715.3215
$ws.Range("J136").Value = 1935.7142
$ws.Range("K136").Value = 8145.9645
$ws.Range("L136").Value = 5807.142599999999
$ws.Range("M136").Value = -5595.9645
$ws.Range("N136").Value = -10907.1426

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9998.75
$ws.Range("J62").Value = 9998.75
$ws.Range("L62").Value = 29996.25
$ws.Range("N62").Value = -31368.25
$ws.Range("H65").Value = 9998.75
$ws.Range("J65").Value = 9998.75
$ws.Range("L65").Value = 89988.75
$ws.Range("N65").Value = -96852.75
$ws.Range("H69").Value = 13500
$ws.Range("I69").Value = 7000
$ws.Range("K69").Value = 21000
$ws.Range("M69").Value = -20189
$ws.Range("H72").Value = 13500
$ws.Range("I72").Value = 7000
$ws.Range("K72").Value = 63000
$ws.Range("M72").Value = -58944
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H113").Value = 1533.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1712.625
$ws.Range("I107").Value = 1283.6666
$ws.Range("K107").Value = 1283.6666
$ws.Range("M107").Value = 636.3334
$ws.Range("H122").Value = 2401.25
$ws.Range("I122").Value = 2401.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7203.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4753.75
$ws.Range("N122").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10517.5
$ws.Range("I22").Value = 12724.333
$ws.Range("K22").Value = 12724.333
$ws.Range("M22").Value = -12429.333
$ws.Range("H27").Value = 10517.5
$ws.Range("I27").Value = 12724.333
$ws.Range("K27").Value = 12724.333
$ws.Range("M27").Value = -12617.333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2526.923
$ws.Range("I132").Value = 1517.1111
$ws.Range("J132").Value = 4799
$ws.Range("K132").Value = 4551.3333
$ws.Range("L132").Value = 14397
$ws.Range("M132").Value = -2021.3333
$ws.Range("N132").Value = -19457
